$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.914.74"
$ws.Range("E2").Value = "  +2.95%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.450.05"
$ws.Range("E3").Value = "  +2.44%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.47"
$ws.Range("E5").Value = "  +4.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "188.88"
$ws.Range("E6").Value = "  +8.33%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.442.58"
$ws.Range("E8").Value = "  +2.52%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  -0.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.646"
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "57.30"
$ws.Range("E12").Value = "  +6.67%  "
$ws.Range("E13").Value = "  -1.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.47"
$ws.Range("E14").Value = "  +3.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.991.63"
$ws.Range("E15").Value = "  +2.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.84"
$ws.Range("E16").Value = "  +3.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.446.93"
$ws.Range("E17").Value = "  +2.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.894.33"
$ws.Range("E18").Value = "  +3.16%  "
$ws.Range("E19").Value = "  +1.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.11"
$ws.Range("E20").Value = "  +1.82%  "
$ws.Range("E21").Value = "  +3.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "480.95"
$ws.Range("E22").Value = "  +5.28%  "
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.12"
$ws.Range("E23").Value = "  +24.30%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.36"
$ws.Range("E24").Value = "  +9.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.35"
$ws.Range("E25").Value = "  +7.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "89.49"
$ws.Range("E26").Value = "  +2.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.99"
$ws.Range("E27").Value = "  +4.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.99"
$ws.Range("E28").Value = "  +2.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.09"
$ws.Range("E29").Value = "  +4.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.28"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.46"
$ws.Range("E31").Value = "  +13.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.81"
$ws.Range("E32").Value = "  +3.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "596.82"
$ws.Range("E33").Value = "  +3.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "64.17"
$ws.Range("E34").Value = "  +2.11%  "
$ws.Range("E35").Value = "  +4.36%  "
$ws.Range("E36").Value = "  +5.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.50"
$ws.Range("E38").Value = "  +5.06%  "
$ws.Range("E39").Value = "  +4.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.51"
$ws.Range("E40").Value = "  -4.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0756"
$ws.Range("E41").Value = "  +2.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.238.50"
$ws.Range("E42").Value = "  +5.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.93"
$ws.Range("E43").Value = "  +5.89%  "
$ws.Range("E44").Value = "  +3.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.83"
$ws.Range("E45").Value = "  +25.91%  "
$ws.Range("E46").Value = "  +3.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.25"
$ws.Range("E47").Value = "  +2.61%  "
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.33"
$ws.Range("E49").Value = "  +15.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.66"
$ws.Range("E51").Value = "  +4.87%  "
